# Commit message: "changing document, table attributes to lowerCamelCase"
#
# The ObjTables header markers written into the first cell(s) of each
# worksheet used upper-camel-case attribute names (ObjTablesVersion, Type,
# Id). They are renamed here to lowerCamelCase (objTablesVersion, type, id)
# while leaving everything else (styling, layout, other cell values)
# untouched.

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item(1)
$wsMain.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$wsMain.Range("A2").Value = "!!ObjTables type='Data' id='MainRoot'"

$wsTests = $wb.Worksheets.Item(2)
$wsTests.Range("A1").Value = "!!ObjTables type='Data' id='Test'"
